$d = $word.ActiveDocument

# --- Update the LTMRdata citation for the updated data publication ---

# 1. Fish survey count grew from 9 to 10
$found1 = $d.Content.Find.Execute(
    "data from 9 fish surveys", $true, $false, $false, $false, $false,
    $true, 1, $false, "data from 10 fish surveys", 2)
Write-Output ("Updated survey count: " + $found1)

# 2. Zenodo DOI now points at the new record (13155050), version text (v2.1.0) unchanged
$found2 = $d.Content.Find.Execute(
    "doi:10.5281/zenodo.6048977", $true, $false, $false, $false, $false,
    $true, 1, $false, "doi:10.5281/zenodo.13155050", 2)
Write-Output ("Updated Zenodo DOI: " + $found2)
